$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "repayment_20250901_20250926"

# Row 2: Debora Retima Sihombing
$ws.Cells.Item(2, 4).Value = 76
$ws.Cells.Item(2, 8).Value = 31.199
$ws.Cells.Item(2, 10).Value = 22
$c = $ws.Cells.Item(2, 5)
$c.NumberFormat = "@"
$c.Value = "51,587,267.00"
$c.Style = "Normal"
$c = $ws.Cells.Item(2, 7)
$c.NumberFormat = "@"
$c.Value = "11.58"
$c.Style = "Normal"
$c = $ws.Cells.Item(2, 11)
$c.NumberFormat = "@"
$c.Value = "8.74"
$c.Style = "Normal"
$c = $ws.Cells.Item(2, 12)
$c.NumberFormat = "@"
$c.Value = "5.31"
$c.Style = "Normal"

# Row 3: Romli
$ws.Cells.Item(3, 4).Value = 69
$ws.Cells.Item(3, 8).Value = 33.759
$c = $ws.Cells.Item(3, 5)
$c.NumberFormat = "@"
$c.Value = "47,627,660.00"
$c.Style = "Normal"
$c = $ws.Cells.Item(3, 7)
$c.NumberFormat = "@"
$c.Value = "10.85"
$c.Style = "Normal"
$c = $ws.Cells.Item(3, 11)
$c.NumberFormat = "@"
$c.Value = "2.73"
$c.Style = "Normal"
$c = $ws.Cells.Item(3, 12)
$c.NumberFormat = "@"
$c.Value = "3.84"
$c.Style = "Normal"

# Row 4: Aldi Taufik
$ws.Cells.Item(4, 4).Value = 59
$ws.Cells.Item(4, 8).Value = 25.214
$c = $ws.Cells.Item(4, 5)
$c.NumberFormat = "@"
$c.Value = "39,587,893.00"
$c.Style = "Normal"
$c = $ws.Cells.Item(4, 7)
$c.NumberFormat = "@"
$c.Value = "8.54"
$c.Style = "Normal"
$c = $ws.Cells.Item(4, 11)
$c.NumberFormat = "@"
$c.Value = "4.83"
$c.Style = "Normal"
$c = $ws.Cells.Item(4, 12)
$c.NumberFormat = "@"
$c.Value = "5.73"
$c.Style = "Normal"

# Row 5: Yandi Nugraha
$ws.Cells.Item(5, 4).Value = 86
$ws.Cells.Item(5, 8).Value = 27.708
$ws.Cells.Item(5, 10).Value = 23
$c = $ws.Cells.Item(5, 5)
$c.NumberFormat = "@"
$c.Value = "49,822,435.00"
$c.Style = "Normal"
$c = $ws.Cells.Item(5, 7)
$c.NumberFormat = "@"
$c.Value = "11.70"
$c.Style = "Normal"
$c = $ws.Cells.Item(5, 11)
$c.NumberFormat = "@"
$c.Value = "5.51"
$c.Style = "Normal"
$c = $ws.Cells.Item(5, 12)
$c.NumberFormat = "@"
$c.Value = "5.97"
$c.Style = "Normal"

# Row 6: Axl Wicaksono
$ws.Cells.Item(6, 4).Value = 62
$ws.Cells.Item(6, 8).Value = 27.216
$c = $ws.Cells.Item(6, 5)
$c.NumberFormat = "@"
$c.Value = "40,700,276.00"
$c.Style = "Normal"
$c = $ws.Cells.Item(6, 7)
$c.NumberFormat = "@"
$c.Value = "9.32"
$c.Style = "Normal"
$c = $ws.Cells.Item(6, 11)
$c.NumberFormat = "@"
$c.Value = "2.87"
$c.Style = "Normal"
$c = $ws.Cells.Item(6, 12)
$c.NumberFormat = "@"
$c.Value = "2.88"
$c.Style = "Normal"

# Row 7: Riska Nurlita
$ws.Cells.Item(7, 4).Value = 91
$ws.Cells.Item(7, 8).Value = 18.6
$ws.Cells.Item(7, 10).Value = 22
$c = $ws.Cells.Item(7, 5)
$c.NumberFormat = "@"
$c.Value = "57,374,493.00"
$c.Style = "Normal"
$c = $ws.Cells.Item(7, 7)
$c.NumberFormat = "@"
$c.Value = "12.23"
$c.Style = "Normal"
$c = $ws.Cells.Item(7, 11)
$c.NumberFormat = "@"
$c.Value = "4.53"
$c.Style = "Normal"
$c = $ws.Cells.Item(7, 12)
$c.NumberFormat = "@"
$c.Value = "5.28"
$c.Style = "Normal"

# Row 8: Annisa Putri Restu
$ws.Cells.Item(8, 4).Value = 75
$ws.Cells.Item(8, 8).Value = 40.223
$c = $ws.Cells.Item(8, 5)
$c.NumberFormat = "@"
$c.Value = "57,713,944.00"
$c.Style = "Normal"
$c = $ws.Cells.Item(8, 7)
$c.NumberFormat = "@"
$c.Value = "13.16"
$c.Style = "Normal"
$c = $ws.Cells.Item(8, 11)
$c.NumberFormat = "@"
$c.Value = "3.47"
$c.Style = "Normal"
$c = $ws.Cells.Item(8, 12)
$c.NumberFormat = "@"
$c.Value = "4.37"
$c.Style = "Normal"

# Row 9: Azizah Rahmawati
$ws.Cells.Item(9, 4).Value = 61
$ws.Cells.Item(9, 8).Value = 18.124
$ws.Cells.Item(9, 10).Value = 13
$c = $ws.Cells.Item(9, 5)
$c.NumberFormat = "@"
$c.Value = "49,569,385.00"
$c.Style = "Normal"
$c = $ws.Cells.Item(9, 7)
$c.NumberFormat = "@"
$c.Value = "10.18"
$c.Style = "Normal"
$c = $ws.Cells.Item(9, 11)
$c.NumberFormat = "@"
$c.Value = "4.19"
$c.Style = "Normal"
$c = $ws.Cells.Item(9, 12)
$c.NumberFormat = "@"
$c.Value = "3.10"
$c.Style = "Normal"

# Row 10: Erlangga Hutama
$ws.Cells.Item(10, 4).Value = 53
$ws.Cells.Item(10, 8).Value = 19.964
$c = $ws.Cells.Item(10, 5)
$c.NumberFormat = "@"
$c.Value = "34,597,319.00"
$c.Style = "Normal"
$c = $ws.Cells.Item(10, 7)
$c.NumberFormat = "@"
$c.Value = "8.70"
$c.Style = "Normal"
$c = $ws.Cells.Item(10, 11)
$c.NumberFormat = "@"
$c.Value = "2.97"
$c.Style = "Normal"
$c = $ws.Cells.Item(10, 12)
$c.NumberFormat = "@"
$c.Value = "2.84"
$c.Style = "Normal"

# Row 11: Erick Ervan Dewanggga
$ws.Cells.Item(11, 4).Value = 64
$ws.Cells.Item(11, 8).Value = 19.267
$c = $ws.Cells.Item(11, 5)
$c.NumberFormat = "@"
$c.Value = "56,061,755.00"
$c.Style = "Normal"
$c = $ws.Cells.Item(11, 7)
$c.NumberFormat = "@"
$c.Value = "12.17"
$c.Style = "Normal"
$c = $ws.Cells.Item(11, 11)
$c.NumberFormat = "@"
$c.Value = "7.68"
$c.Style = "Normal"
$c = $ws.Cells.Item(11, 12)
$c.NumberFormat = "@"
$c.Value = "6.04"
$c.Style = "Normal"

# Row 12: Ridhoi Berkat Zebua
$ws.Cells.Item(12, 4).Value = 71
$ws.Cells.Item(12, 8).Value = 36.587
$c = $ws.Cells.Item(12, 5)
$c.NumberFormat = "@"
$c.Value = "43,601,161.00"
$c.Style = "Normal"
$c = $ws.Cells.Item(12, 7)
$c.NumberFormat = "@"
$c.Value = "9.88"
$c.Style = "Normal"
$c = $ws.Cells.Item(12, 11)
$c.NumberFormat = "@"
$c.Value = "5.38"
$c.Style = "Normal"
$c = $ws.Cells.Item(12, 12)
$c.NumberFormat = "@"
$c.Value = "5.73"
$c.Style = "Normal"

# Row 13: Fadilah Damayanti
$ws.Cells.Item(13, 4).Value = 61
$ws.Cells.Item(13, 8).Value = 23.716
$c = $ws.Cells.Item(13, 5)
$c.NumberFormat = "@"
$c.Value = "40,816,725.00"
$c.Style = "Normal"
$c = $ws.Cells.Item(13, 7)
$c.NumberFormat = "@"
$c.Value = "9.08"
$c.Style = "Normal"
$c = $ws.Cells.Item(13, 11)
$c.NumberFormat = "@"
$c.Value = "3.17"
$c.Style = "Normal"
$c = $ws.Cells.Item(13, 12)
$c.NumberFormat = "@"
$c.Value = "2.68"
$c.Style = "Normal"

# Row 14: Nur Halim
$ws.Cells.Item(14, 4).Value = 84
$ws.Cells.Item(14, 8).Value = 14.457
$c = $ws.Cells.Item(14, 5)
$c.NumberFormat = "@"
$c.Value = "69,983,773.00"
$c.Style = "Normal"
$c = $ws.Cells.Item(14, 7)
$c.NumberFormat = "@"
$c.Value = "15.61"
$c.Style = "Normal"
$c = $ws.Cells.Item(14, 11)
$c.NumberFormat = "@"
$c.Value = "8.38"
$c.Style = "Normal"
$c = $ws.Cells.Item(14, 12)
$c.NumberFormat = "@"
$c.Value = "3.86"
$c.Style = "Normal"

# Row 15: Adistira Winditya P
$ws.Cells.Item(15, 4).Value = 61
$ws.Cells.Item(15, 8).Value = 20.984
$c = $ws.Cells.Item(15, 5)
$c.NumberFormat = "@"
$c.Value = "38,423,831.00"
$c.Style = "Normal"
$c = $ws.Cells.Item(15, 7)
$c.NumberFormat = "@"
$c.Value = "8.55"
$c.Style = "Normal"
$c = $ws.Cells.Item(15, 11)
$c.NumberFormat = "@"
$c.Value = "3.22"
$c.Style = "Normal"
$c = $ws.Cells.Item(15, 12)
$c.NumberFormat = "@"
$c.Value = "3.37"
$c.Style = "Normal"

# Row 16: Sucika Wardani
$ws.Cells.Item(16, 4).Value = 65
$ws.Cells.Item(16, 8).Value = 15.092
$ws.Cells.Item(16, 10).Value = 15
$c = $ws.Cells.Item(16, 5)
$c.NumberFormat = "@"
$c.Value = "41,773,083.00"
$c.Style = "Normal"
$c = $ws.Cells.Item(16, 7)
$c.NumberFormat = "@"
$c.Value = "9.23"
$c.Style = "Normal"
$c = $ws.Cells.Item(16, 11)
$c.NumberFormat = "@"
$c.Value = "2.65"
$c.Style = "Normal"
$c = $ws.Cells.Item(16, 12)
$c.NumberFormat = "@"
$c.Value = "3.57"
$c.Style = "Normal"

# Row 17: Wasti Feronika Sihombing
$ws.Cells.Item(17, 4).Value = 57
$ws.Cells.Item(17, 8).Value = 29.507
$c = $ws.Cells.Item(17, 5)
$c.NumberFormat = "@"
$c.Value = "44,485,160.00"
$c.Style = "Normal"
$c = $ws.Cells.Item(17, 7)
$c.NumberFormat = "@"
$c.Value = "10.12"
$c.Style = "Normal"
$c = $ws.Cells.Item(17, 11)
$c.NumberFormat = "@"
$c.Value = "3.17"
$c.Style = "Normal"
$c = $ws.Cells.Item(17, 12)
$c.NumberFormat = "@"
$c.Value = "4.35"
$c.Style = "Normal"

# Row 18: Nuraini
$ws.Cells.Item(18, 4).Value = 50
$ws.Cells.Item(18, 8).Value = 19.624
$c = $ws.Cells.Item(18, 5)
$c.NumberFormat = "@"
$c.Value = "38,229,259.00"
$c.Style = "Normal"
$c = $ws.Cells.Item(18, 7)
$c.NumberFormat = "@"
$c.Value = "9.89"
$c.Style = "Normal"
$c = $ws.Cells.Item(18, 11)
$c.NumberFormat = "@"
$c.Value = "3.63"
$c.Style = "Normal"
$c = $ws.Cells.Item(18, 12)
$c.NumberFormat = "@"
$c.Value = "2.83"
$c.Style = "Normal"

